$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table 1 (rows 3-7, 9, 13-14): columns B (d0c_glmm) and C (d0c_jm)
$ws.Range("B3").Value = -1.0099514125774056
$ws.Range("C3").Value = -0.99628184150496013

$ws.Range("B4").Value = -1.0020966887339296
$ws.Range("C4").Value = -0.98534907595284949

$ws.Range("B5").Value = -0.98231031340247543
$ws.Range("C5").Value = -0.96371712264362719

$ws.Range("B6").Value = -1.010322508677292
$ws.Range("C6").Value = -0.9901294300598823

$ws.Range("B7").Value = -0.9500542817907317
$ws.Range("C7").Value = -0.92876948679706606

$ws.Range("B9").Value = 0.070506206282825126
$ws.Range("C9").Value = 0.069869312552049354

$ws.Range("B13").Value = 0.12066991499573686
$ws.Range("C13").Value = 0.11930353682840983

$ws.Range("B14").Value = 1.7084690953288504
$ws.Range("C14").Value = 1.7010545977157492

# Table 2 (rows 16-20): column C only
$ws.Range("C16").Value = -0.2232463549633503
$ws.Range("C17").Value = -0.084037505274748722
$ws.Range("C18").Value = 0.029509486888930313
$ws.Range("C19").Value = -1.5087899445940627
$ws.Range("C20").Value = 0.025943025357203889

# Row 21: Number of observations
$ws.Range("B21").Value = 8267
$ws.Range("C21").Value = 10204
